$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Persona Física"
$ws.Range("B4").Value = "Tercer atributo"
